# Merge the two runs "Home Owner" + " Flag (H = Homeowner; U = Unknown)"
# (previously split across a w:proofErr gramStart/gramEnd pair) into a
# single run reading "Home Owner Flag (H = Homeowner; U = Unknown)".
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Home Owner Flag (H = Homeowner; U = Unknown)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Home Owner Flag (H = Homeowner; U = Unknown)",
    2)
